# Auto-generated edit script: updates market-price-derived profit columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the scheduled price refresh.
$wb = $excel.ActiveWorkbook

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2180.5
$ws.Range("I62").Value = 2180.5
$ws.Range("K62").Value = 2180.5
$ws.Range("M62").Value = -1556.5

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2180.5
$ws.Range("I65").Value = 2180.5
$ws.Range("K65").Value = 10902.5
$ws.Range("M65").Value = -7782.5

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 714.9
$ws.Range("I80").Value = 749.6667
$ws.Range("J80").Value = 700
$ws.Range("K80").Value = 2249.0001
$ws.Range("L80").Value = 2100
$ws.Range("M80").Value = -1251.0001
$ws.Range("N80").Value = -4096

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 714.9
$ws.Range("I83").Value = 749.6667
$ws.Range("J83").Value = 700
$ws.Range("K83").Value = 6747.0003
$ws.Range("L83").Value = 6300
$ws.Range("M83").Value = -1755.0003
$ws.Range("N83").Value = -16284

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3108.4
$ws.Range("I86").Value = 2198.5
$ws.Range("K86").Value = 2198.5
$ws.Range("M86").Value = -1075.5

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 3489.5
$ws.Range("I88").Value = 2627.8572
$ws.Range("K88").Value = 2627.8572
$ws.Range("M88").Value = -2221.8572

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3108.4
$ws.Range("I89").Value = 2198.5
$ws.Range("K89").Value = 10992.5
$ws.Range("M89").Value = -5376.5

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 3489.5
$ws.Range("I91").Value = 2627.8572
$ws.Range("K91").Value = 2627.8572
$ws.Range("M91").Value = -1223.8572

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1269.5555
$ws.Range("I100").Value = 915.6842
$ws.Range("J100").Value = 2110
$ws.Range("K100").Value = 915.6842
$ws.Range("L100").Value = 2110
$ws.Range("M100").Value = -374.6842
$ws.Range("N100").Value = -3192

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2459.6875
$ws.Range("I106").Value = 2477
$ws.Range("K106").Value = 2477
$ws.Range("M106").Value = -1846

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3974.889
$ws.Range("I113").Value = 3955
$ws.Range("K113").Value = 3955
$ws.Range("M113").Value = -701

# ALC row 136
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 199999.67
$ws.Range("J136").Value = 199999.67
$ws.Range("L136").Value = 199999.67
$ws.Range("N136").Value = -210199.67

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1923.6957
$ws.Range("J138").Value = 2469.3713
$ws.Range("L138").Value = 7408.113899999999
$ws.Range("N138").Value = -17688.1139

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6784.5386
$ws.Range("I61").Value = 5381.758
$ws.Range("K61").Value = 5381.758
$ws.Range("M61").Value = -5169.758

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1867.6154
$ws.Range("I88").Value = 1496.75
$ws.Range("J88").Value = 2032.4445
$ws.Range("K88").Value = 1496.75
$ws.Range("L88").Value = 2032.4445
$ws.Range("M88").Value = -1090.75
$ws.Range("N88").Value = -2844.4445

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1867.6154
$ws.Range("I91").Value = 1496.75
$ws.Range("J91").Value = 2032.4445
$ws.Range("K91").Value = 1496.75
$ws.Range("L91").Value = 2032.4445
$ws.Range("M91").Value = -92.75
$ws.Range("N91").Value = -4840.4445

# ARM row 114
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H114").Value = 30398
$ws.Range("J114").Value = 30398
$ws.Range("L114").Value = 30398
$ws.Range("N114").Value = -39076

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 9299.125
$ws.Range("I132").Value = 7537.5
$ws.Range("K132").Value = 22612.5
$ws.Range("M132").Value = -20082.5

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6784.5386
$ws.Range("I136").Value = 5381.758
$ws.Range("K136").Value = 16145.274
$ws.Range("M136").Value = -13595.274

# ARM row 138
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("N138").ClearContents()

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 675.93335
$ws.Range("I22").Value = 715.6429000000001
$ws.Range("K22").Value = 715.6429000000001
$ws.Range("M22").Value = -542.6429000000001

# BSM row 35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 42611.5
$ws.Range("J35").Value = 45973.8
$ws.Range("L35").Value = 45973.8
$ws.Range("N35").Value = -46593.8

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1549.7333
$ws.Range("I94").Value = 1451.826
$ws.Range("K94").Value = 1451.826
$ws.Range("M94").Value = -1000.826

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4666.5
$ws.Range("I107").Value = 3875
$ws.Range("K107").Value = 3875
$ws.Range("M107").Value = -1955

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4965.8945
$ws.Range("I122").Value = 4917.6665
$ws.Range("J122").Value = 5009.3
$ws.Range("K122").Value = 14752.9995
$ws.Range("L122").Value = 15027.9
$ws.Range("M122").Value = -12302.9995
$ws.Range("N122").Value = -19927.9

# CUL row 82
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 7522500
$ws.Range("J82").Value = 45000
$ws.Range("L82").Value = 135000
$ws.Range("N82").Value = -135812

# CUL row 85
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 7522500
$ws.Range("J85").Value = 45000
$ws.Range("L85").Value = 135000
$ws.Range("N85").Value = -137808

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 24603936
$ws.Range("I129").Value = 47763584
$ws.Range("K129").Value = 143290752
$ws.Range("M129").Value = -143285752

# GSM row 43
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 5771.8184
$ws.Range("I97").Value = 1242.1428
$ws.Range("K97").Value = 1242.1428
$ws.Range("M97").Value = -746.1428000000001

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 287287.94
$ws.Range("I113").Value = 334770.25
$ws.Range("K113").Value = 334770.25
$ws.Range("M113").Value = -332600.25

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3398
$ws.Range("I122").Value = 3590.25
$ws.Range("J122").Value = 3013.5
$ws.Range("K122").Value = 10770.75
$ws.Range("L122").Value = 9040.5
$ws.Range("M122").Value = -8320.75
$ws.Range("N122").Value = -13940.5

# GSM row 124
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 79780
$ws.Range("J124").Value = 79780
$ws.Range("L124").Value = 79780
$ws.Range("N124").Value = -89600

# GSM row 133
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("N133").ClearContents()

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2639.6
$ws.Range("I40").Value = 2456.7856
$ws.Range("K40").Value = 2456.7856
$ws.Range("M40").Value = -2320.7856

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 10792.465
$ws.Range("I46").Value = 5150.8184
$ws.Range("K46").Value = 5150.8184
$ws.Range("M46").Value = -4962.8184

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3316.5
$ws.Range("I122").Value = 2876.3333
$ws.Range("K122").Value = 8628.999899999999
$ws.Range("M122").Value = -6178.999899999999

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9730.65
$ws.Range("I132").Value = 10017.454
$ws.Range("K132").Value = 30052.362
$ws.Range("M132").Value = -27522.362

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1156.4286
$ws.Range("J81").Value = 1499.3334
$ws.Range("L81").Value = 2998.6668
$ws.Range("N81").Value = -5120.6668

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1156.4286
$ws.Range("J84").Value = 1499.3334
$ws.Range("L84").Value = 14993.334
$ws.Range("N84").Value = -25601.334

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4335.4287
$ws.Range("I122").Value = 4157.45
$ws.Range("K122").Value = 12472.35
$ws.Range("M122").Value = -10022.35

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5670.5
$ws.Range("I126").Value = 5226.778
$ws.Range("K126").Value = 15680.334
$ws.Range("M126").Value = -13210.334

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3066.7368
$ws.Range("I132").Value = 2994.879
$ws.Range("J132").Value = 3541
$ws.Range("K132").Value = 8984.636999999999
$ws.Range("L132").Value = 10623
$ws.Range("M132").Value = -6454.636999999999
$ws.Range("N132").Value = -15683

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3957.4546
$ws.Range("I136").Value = 2434.4482
$ws.Range("J136").Value = 14999.25
$ws.Range("K136").Value = 7303.344599999999
$ws.Range("L136").Value = 44997.75
$ws.Range("M136").Value = -4753.344599999999
$ws.Range("N136").Value = -50097.75
